$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCells = @("D5", "D6", "D10", "D11", "D12", "D19", "D22", "D23", "D27", "D28", "D30", "D40", "D41", "D44", "D45", "D49", "D51")
foreach ($c in $numCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.058.68"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "2.304.88"

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "306.17"
$ws.Range("E5").Value = "  +2.65%  "

$ws.Range("D6").Value = "97.56"
$ws.Range("E6").Value = "  +1.36%  "

$ws.Range("E7").Value = "  -0.94%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").Value = "35.56"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").Value = "18.72"
$ws.Range("E12").Value = "  +6.52%  "

$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("E14").Value = "  +2.29%  "

$ws.Range("D15").Value = "2.663.69"

$ws.Range("D16").Value = "2.299.86"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").Value = "42.942.62"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  +2.10%  "

$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").Value = "67.38"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").Value = "236.71"
$ws.Range("E23").Value = "  -1.62%  "

$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("E26").Value = "  +0.67%  "

$ws.Range("D27").Value = "24.94"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("D28").Value = "166.89"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("D30").Value = "9.06"
$ws.Range("E30").Value = "  +0.72%  "

$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("E33").Value = "  +6.48%  "

$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("E35").Value = "  -5.21%  "

$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").Value = "1.996.90"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").Value = "10.37"
$ws.Range("E44").Value = "  +3.03%  "

$ws.Range("D45").Value = "17.98"
$ws.Range("E45").Value = "  +6.05%  "

$ws.Range("E46").Value = "  +1.99%  "

$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("E48").Value = "  +4.51%  "

$ws.Range("D49").Value = "53.68"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").Value = "2.531.58"
$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "71.59"
$ws.Range("E51").Value = "  -0.20%  "

foreach ($c in $numCells) {
    $ws.Range($c).Style = "Normal"
}
